$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 9; $r++) {
    $cellC = $ws.Range("C$r")
    $cellC.NumberFormat = "@"
    $cellC.Value = "2024-01-02"
    $cellC.Style = "常规"

    $cellH = $ws.Range("H$r")
    $cellH.NumberFormat = "@"
    $cellH.Value = "16"
    $cellH.Style = "常规"

    $cellI = $ws.Range("I$r")
    $cellI.NumberFormat = "@"
    $cellI.Value = "15"
    $cellI.Style = "常规"
}
